$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New symbolValWins data (rows 2-21), reordered/updated with collectorNum now
# encoded into the symbol value (symbol = base*100 + collectorNum).
$data = @(
  @(901, 16, 15, 45, 60, 60),
  @(902, 1, 0, 0, 0, 0),
  @(301, 6, 45, 30, 60, 45),
  @(401, 9, 48, 67, 75, 45),
  @(701, 3, 90, 45, 97, 15),
  @(601, 9, 60, 67, 60, 42),
  @(101, 9, 30, 15, 60, 15),
  @(1001, 18, 30, 75, 60, 72),
  @(501, 9, 52, 30, 75, 45),
  @(1201, 2, 10, 10, 10, 10),
  @(201, 9, 30, 15, 45, 30),
  @(801, 3, 67, 65, 52, 45),
  @(1202, 2, 10, 10, 10, 10),
  @(1203, 3, 15, 15, 15, 15),
  @(802, 0, 4, 5, 4, 0),
  @(3, 0, 3, 3, 3, 3),
  @(2, 0, 2, 2, 2, 2),
  @(1, 0, 2, 2, 2, 2),
  @(502, 0, 4, 0, 0, 0),
  @(1101, 0, 15, 30, 30, 0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $data[$i]
  for ($j = 0; $j -lt $row.Length; $j++) {
    $ws.Cells.Item($i + 2, $j + 1).Value = $row[$j]
  }
}
